# Auto-generated edit script: updates computed market-price / profit columns (H:N)
# across several sheets, per the scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1887737.6
$ws.Range("J17").Value = 1887737.6
$ws.Range("L17").Value = 5663212.800000001
$ws.Range("N17").Value = -5663548.800000001

$ws.Range("H135").Value = 1255.7322
$ws.Range("I135").Value = 752.0833
$ws.Range("J135").Value = 4277.625
$ws.Range("K135").Value = 6768.7497
$ws.Range("L135").Value = 38498.625
$ws.Range("M135").Value = -4233.7497
$ws.Range("N135").Value = -43568.625

$ws.Range("H137").Value = 1011.3038
$ws.Range("I137").Value = 827.5069999999999
$ws.Range("J137").Value = 2642.5
$ws.Range("K137").Value = 2482.521
$ws.Range("L137").Value = 7927.5
$ws.Range("M137").Value = 67.47900000000027
$ws.Range("N137").Value = -13027.5

$ws.Range("H138").Value = 3024.963
$ws.Range("I138").Value = 1384
$ws.Range("J138").Value = 4928.48
$ws.Range("K138").Value = 4152
$ws.Range("L138").Value = 14785.44
$ws.Range("M138").Value = 988
$ws.Range("N138").Value = -25065.44

$ws.Range("H141").Value = 1702.2565
$ws.Range("I141").Value = 1248.8
$ws.Range("J141").Value = 3213.7778
$ws.Range("K141").Value = 3746.4
$ws.Range("L141").Value = 9641.3334
$ws.Range("M141").Value = 1433.6
$ws.Range("N141").Value = -20001.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3248.9
$ws.Range("I32").Value = 2295.8489
$ws.Range("J32").Value = 9103.357
$ws.Range("K32").Value = 2295.8489
$ws.Range("L32").Value = 9103.357
$ws.Range("M32").Value = -2008.8489
$ws.Range("N32").Value = -9677.357

$ws.Range("H61").Value = 3332.2449
$ws.Range("I61").Value = 4164.606
$ws.Range("J61").Value = 1615.5
$ws.Range("K61").Value = 4164.606
$ws.Range("L61").Value = 1615.5
$ws.Range("M61").Value = -3952.606
$ws.Range("N61").Value = -2039.5

$ws.Range("H88").Value = 1863.25
$ws.Range("I88").Value = 1843.7142
$ws.Range("J88").Value = 2000
$ws.Range("K88").Value = 1843.7142
$ws.Range("L88").Value = 2000
$ws.Range("M88").Value = -1437.7142
$ws.Range("N88").Value = -2812

$ws.Range("H91").Value = 1863.25
$ws.Range("I91").Value = 1843.7142
$ws.Range("J91").Value = 2000
$ws.Range("K91").Value = 1843.7142
$ws.Range("L91").Value = 2000
$ws.Range("M91").Value = -439.7141999999999
$ws.Range("N91").Value = -4808

$ws.Range("H136").Value = 3332.2449
$ws.Range("I136").Value = 4164.606
$ws.Range("J136").Value = 1615.5
$ws.Range("K136").Value = 12493.818
$ws.Range("L136").Value = 4846.5
$ws.Range("M136").Value = -9943.817999999999
$ws.Range("N136").Value = -9946.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877

$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384

$ws.Range("H99").Value = 71429680
$ws.Range("I99").Value = 83334090
$ws.Range("J99").Value = 3250
$ws.Range("K99").Value = 83334090
$ws.Range("L99").Value = 3250
$ws.Range("M99").Value = -83332592
$ws.Range("N99").Value = -6246

$ws.Range("H134").Value = 2579.6956
$ws.Range("I134").Value = 2532.242
$ws.Range("K134").Value = 7596.726000000001
$ws.Range("M134").Value = -5061.726000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5663.1055
$ws.Range("I31").Value = 1454.3478
$ws.Range("K31").Value = 1454.3478
$ws.Range("M31").Value = -1159.3478

$ws.Range("H34").Value = 5663.1055
$ws.Range("I34").Value = 1454.3478
$ws.Range("K34").Value = 1454.3478
$ws.Range("M34").Value = -1252.3478

$ws.Range("H58").Value = 1032.742
$ws.Range("I58").Value = 646.6977000000001
$ws.Range("J58").Value = 1906.421
$ws.Range("K58").Value = 646.6977000000001
$ws.Range("L58").Value = 1906.421
$ws.Range("M58").Value = -443.6977000000001
$ws.Range("N58").Value = -2312.421

$ws.Range("H62").Value = 6130.8335
$ws.Range("I62").Value = 6557
$ws.Range("K62").Value = 6557
$ws.Range("M62").Value = -5933

$ws.Range("H65").Value = 6130.8335
$ws.Range("I65").Value = 6557
$ws.Range("K65").Value = 32785
$ws.Range("M65").Value = -29665

$ws.Range("H99").Value = 13903200
$ws.Range("I99").Value = 17180
$ws.Range("J99").Value = 31260724
$ws.Range("K99").Value = 17180
$ws.Range("L99").Value = 31260724
$ws.Range("M99").Value = -15682
$ws.Range("N99").Value = -31263720

$ws.Range("H126").Value = 13903200
$ws.Range("I126").Value = 17180
$ws.Range("J126").Value = 31260724
$ws.Range("K126").Value = 51540
$ws.Range("L126").Value = 93782172
$ws.Range("M126").Value = -49070
$ws.Range("N126").Value = -93787112

$ws.Range("H132").Value = 1258.6621
$ws.Range("I132").Value = 906.5848999999999
$ws.Range("J132").Value = 2147.238
$ws.Range("K132").Value = 2719.7547
$ws.Range("L132").Value = 6441.714
$ws.Range("M132").Value = -189.7547
$ws.Range("N132").Value = -11501.714

$ws.Range("H134").Value = 1470.8
$ws.Range("I134").Value = 1573.2106
$ws.Range("J134").Value = 1262.3214
$ws.Range("K134").Value = 4719.6318
$ws.Range("L134").Value = 3786.9642
$ws.Range("M134").Value = -2184.6318
$ws.Range("N134").Value = -8856.9642

$ws.Range("H136").Value = 1032.742
$ws.Range("I136").Value = 646.6977000000001
$ws.Range("J136").Value = 1906.421
$ws.Range("K136").Value = 1940.0931
$ws.Range("L136").Value = 5719.263
$ws.Range("M136").Value = 609.9069
$ws.Range("N136").Value = -10819.263

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 97134.61
$ws.Range("I5").Value = 178.89743
$ws.Range("J5").Value = 261537.78
$ws.Range("K5").Value = 536.6922900000001
$ws.Range("L5").Value = 784613.34
$ws.Range("M5").Value = -424.6922900000001
$ws.Range("N5").Value = -784837.34

$ws.Range("H12").Value = 3448393.2
$ws.Range("I12").Value = 8333421.5
$ws.Range("K12").Value = 25000264.5
$ws.Range("M12").Value = -25000091.5

$ws.Range("H113").Value = 345395.28
$ws.Range("I113").Value = 632.9167
$ws.Range("J113").Value = 588756.9399999999
$ws.Range("K113").Value = 1898.7501
$ws.Range("L113").Value = 1766270.82
$ws.Range("M113").Value = 271.2499
$ws.Range("N113").Value = -1770610.82

$ws.Range("H116").Value = 500
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 500
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 1500
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -8384

$ws.Range("H131").Value = 1786623.1
$ws.Range("I131").Value = 6250659.5
$ws.Range("J131").Value = 1008.65
$ws.Range("K131").Value = 18751978.5
$ws.Range("L131").Value = 3025.95
$ws.Range("M131").Value = -18746938.5
$ws.Range("N131").Value = -13105.95

$ws.Range("H134").Value = 6587475.5
$ws.Range("I134").Value = 14713363
$ws.Range("J134").Value = 9376.190000000001
$ws.Range("K134").Value = 44140089
$ws.Range("L134").Value = 28128.57
$ws.Range("M134").Value = -44135019
$ws.Range("N134").Value = -38268.57

$ws.Range("H135").Value = 97134.61
$ws.Range("I135").Value = 178.89743
$ws.Range("J135").Value = 261537.78
$ws.Range("K135").Value = 1610.07687
$ws.Range("L135").Value = 2353840.02
$ws.Range("M135").Value = 924.9231299999999
$ws.Range("N135").Value = -2358910.02

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1326.0704
$ws.Range("I132").Value = 969.8909
$ws.Range("J132").Value = 2550.4375
$ws.Range("K132").Value = 2909.6727
$ws.Range("L132").Value = 7651.3125
$ws.Range("M132").Value = -379.6727000000001
$ws.Range("N132").Value = -12711.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 427.7143
$ws.Range("I16").Value = 427.7143
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 427.7143
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -257.7143
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 898.4792
$ws.Range("I132").Value = 702.70966
$ws.Range("J132").Value = 1255.4706
$ws.Range("K132").Value = 2108.12898
$ws.Range("L132").Value = 3766.4118
$ws.Range("M132").Value = 421.87102
$ws.Range("N132").Value = -8826.4118

$ws.Range("H136").Value = 5954197
$ws.Range("I136").Value = 1993.4728
$ws.Range("J136").Value = 17242858
$ws.Range("K136").Value = 5980.4184
$ws.Range("L136").Value = 51728574
$ws.Range("M136").Value = -3430.4184
$ws.Range("N136").Value = -51733674
